$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the diff. Columns B/C/E are plain text (safe to set directly).
# Column D values are text-formatted numbers (e.g. "1.002", "0.00001029", "29.458.70")
# that Excel would otherwise auto-convert to numeric/scientific values, so we force
# the cell to Text format before assigning the literal string.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.458.70"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.909.12"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.47"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4663"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4078"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.75"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08020"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.30"
$ws.Range("E12").Value = "  +2.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.943.36"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.944"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.134"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.14"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001029"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.462.04"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.538"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.48"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.211"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.111.24"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.62"
$ws.Range("E27").Value = "  -2.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.74"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.713"
$ws.Range("E29").Value = "  +5.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.124"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.82"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.072"
$ws.Range("E32").Value = "  +9.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09443"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.419"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.577"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.391"
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06080"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.383"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.171"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5867"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.301"
$ws.Range("E44").Value = "  +3.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07753"
$ws.Range("E45").Value = "  +10.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.374"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5538"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.09"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.28"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2933"
$ws.Range("E51").Value = "  +3.85%  "
